$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Conditioning column (H) rebalancing ---
# H4 becomes a plain entered value (was 15, now 23).
$ws.Range("H4").Value = 23

# H5:H12 are filled down with the same relative formula, each cell
# referencing the Conditioning value one row above it.
$ws.Range("H5:H12").FormulaR1C1 = "=ROUND(1.3125*R[-1]C,0)"

# --- Selection moves to where the edit left off ---
$ws.Range("H12").Select()
